# Apply "想去人数" (interested-count) bumps across the three sheets that
# carry per-event attendee numbers, plus a newly-scraped "演出" (Performance)
# row describing a freshly-listed event.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) updates ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 80
$ws1.Range("F3").Value = 173
$ws1.Range("F8").Value = 1640
$ws1.Range("G8").Value = 60
$ws1.Range("F9").Value = 7457
$ws1.Range("F11").Value = 7651
$ws1.Range("F12").Value = 18
$ws1.Range("F15").Value = 6177
$ws1.Range("F16").Value = 3264
$ws1.Range("F17").Value = 3628
$ws1.Range("F18").Value = 19
$ws1.Range("F19").Value = 12
$ws1.Range("F21").Value = 30
$ws1.Range("F22").Value = 446
$ws1.Range("F24").Value = 285
$ws1.Range("F26").Value = 3632
$ws1.Range("F30").Value = 259
$ws1.Range("F31").Value = 1097
$ws1.Range("F32").Value = 64
$ws1.Range("F34").Value = 2621
$ws1.Range("F35").Value = 1464
$ws1.Range("F38").Value = 23
$ws1.Range("F39").Value = 3274
$ws1.Range("F41").Value = 244
$ws1.Range("F44").Value = 481
$ws1.Range("F45").Value = 1284
$ws1.Range("F46").Value = 226
$ws1.Range("F48").Value = 593

# --- Sheet "演出" (Performance) updates ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 235
$ws2.Range("F16").Value = 79

# New event row appended to "演出": copy the formatting of the last data
# row (row 18, which carries the bold/boxed index style on column A) down
# to the new row 19, then fill in the values so the style index is
# preserved instead of left blank.
$srcIndexCell = $ws2.Cells.Item(18, 1)
$dstIndexCell = $ws2.Cells.Item(19, 1)
$srcIndexCell.Copy($dstIndexCell)

$ws2.Cells.Item(19, 1).Value = 18

# Column B holds a plain "yyyy-MM-dd" label, not a real date, in every
# other row (inlineStr). Force text format so assigning the
# date-shaped string doesn't get auto-coerced into a date serial
# number/date-formatted cell, then drop back to the default style so no
# stray numFmt sticks to the cell.
$dateCell = $ws2.Cells.Item(19, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-08-24"
$dateCell.Style = "Normal"

$ws2.Cells.Item(19, 3).Value = "北京·最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会"
$ws2.Cells.Item(19, 4).Value = "亮马桥路40号(近好运街) 北京世纪剧院"
$ws2.Cells.Item(19, 5).Value = "2024.08.24 19:30-08.24 21:00"
$ws2.Cells.Item(19, 6).Value = 0
$ws2.Cells.Item(19, 7).Value = 144
$ws2.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86217"
$ws2.Cells.Item(19, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/BDyblKrJ1716427731729.jpeg"

# --- Sheet "全部类型" (All types) updates ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 80
$ws4.Range("F4").Value = 173
$ws4.Range("F9").Value = 1640
$ws4.Range("G9").Value = 60
$ws4.Range("F10").Value = 235
$ws4.Range("F12").Value = 7457
$ws4.Range("F13").Value = 7651
$ws4.Range("F14").Value = 18
$ws4.Range("F15").Value = 6177
$ws4.Range("F16").Value = 3264
$ws4.Range("F17").Value = 3628
$ws4.Range("F18").Value = 19
$ws4.Range("F19").Value = 12
$ws4.Range("F21").Value = 30
$ws4.Range("F22").Value = 446
$ws4.Range("F23").Value = 285
$ws4.Range("F27").Value = 3632
$ws4.Range("F32").Value = 259
$ws4.Range("F33").Value = 64
$ws4.Range("F35").Value = 2621
$ws4.Range("F36").Value = 1464
$ws4.Range("F39").Value = 79
$ws4.Range("F40").Value = 3274
$ws4.Range("F42").Value = 244
$ws4.Range("F46").Value = 481
$ws4.Range("F47").Value = 1284
$ws4.Range("F48").Value = 226

Write-Host "Edit script completed."
